$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'246.93"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'0.74%"
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(3,4).Value = "'26.42"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'5.01%"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(4,4).Value = "'5.075"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'0.29%"
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(5,4).Value = "'0.05600"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'-0.24%"
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'0.60%"
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(8,4).Value = "'0.8454"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "'0.00%"
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(9,2).Value = "One"
$ws.Cells.Item(9,3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(9,4).Value = "'0.009930"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "'1,561.05%"
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(10,2).Value = "MandalaExchangeToken"
$ws.Cells.Item(10,3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(10,4).Value = "'0.06986"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'0.28%"
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(11,2).Value = "BitrueCoin"
$ws.Cells.Item(11,3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(11,4).Value = "'0.02850"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'0.36%"
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(12,2).Value = "BitMartToken"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(12,4).Value = "'0.09393"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "'-0.22%"
$ws.Cells.Item(12,5).Style = "Normal"
$ws.Cells.Item(13,2).Value = "BitForexToken"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(13,4).Value = "'0.001509"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "'-1.03%"
$ws.Cells.Item(13,5).Style = "Normal"
$ws.Cells.Item(14,4).Value = "'0.006147"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "'-0.34%"
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(15,4).Value = "'3.601"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "'2.86%"
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(16,5).Value = "'0.31%"
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(17,5).Value = "'-1.71%"
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(18,4).Value = "'0.3156"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "'-0.84%"
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(19,5).Value = "'-0.13%"
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(20,4).Value = "'0.03188"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "'-1.90%"
$ws.Cells.Item(20,5).Style = "Normal"
$ws.Cells.Item(21,5).Value = "'0.49%"
$ws.Cells.Item(21,5).Style = "Normal"
$ws.Cells.Item(22,4).Value = "'3.740"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "'0.10%"
$ws.Cells.Item(22,5).Style = "Normal"
$ws.Cells.Item(23,4).Value = "'0.04631"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "'-0.97%"
$ws.Cells.Item(23,5).Style = "Normal"
$ws.Cells.Item(24,5).Value = "'-1.46%"
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(25,5).Value = "'0.19%"
$ws.Cells.Item(25,5).Style = "Normal"
$ws.Cells.Item(26,4).Value = "'0.004586"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "'1.23%"
$ws.Cells.Item(26,5).Style = "Normal"
$ws.Cells.Item(27,4).Value = "'0.00009597"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "'-1.02%"
$ws.Cells.Item(27,5).Style = "Normal"
$ws.Cells.Item(28,4).Value = "'0.0001397"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = "'1.66%"
$ws.Cells.Item(28,5).Style = "Normal"
$ws.Cells.Item(40,4).Value = "'0.03668"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "'0.71%"
$ws.Cells.Item(40,5).Style = "Normal"
$ws.Cells.Item(41,2).Value = "KickToken"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(41,4).Value = "'0.006128"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "'-1.76%"
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(42,2).Value = "BKEXToken"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(42,4).Value = "'0.1057"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "'-21.63%"
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(43,4).Value = "'0.002444"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "'-10.25%"
$ws.Cells.Item(43,5).Style = "Normal"
$ws.Cells.Item(44,4).Value = "'0.008934"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "'10.77%"
$ws.Cells.Item(44,5).Style = "Normal"
$ws.Cells.Item(45,4).Value = "'0.00005272"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "'-0.14%"
$ws.Cells.Item(45,5).Style = "Normal"
$ws.Cells.Item(46,5).Value = "'0.00%"
$ws.Cells.Item(46,5).Style = "Normal"
$ws.Cells.Item(47,5).Value = "'-38.89%"
$ws.Cells.Item(47,5).Style = "Normal"
$ws.Cells.Item(48,4).Value = "'0.002626"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "'28.47%"
$ws.Cells.Item(48,5).Style = "Normal"
$ws.Cells.Item(49,4).Value = "'0.00002099"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "'0.00%"
$ws.Cells.Item(49,5).Style = "Normal"
$ws.Cells.Item(50,4).Value = "'0.0001999"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "'0.00%"
$ws.Cells.Item(50,5).Style = "Normal"
